# Filmvorführung Buchspazierer 20.12.24 Eintritte und Kiosk
# Append three new rows of kiosk sales data for the 20.12.2024 screening.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Date for the new entries: 20.12.2024 (serial 45646)
$newDate = Get-Date -Year 2024 -Month 12 -Day 20 -Hour 0 -Minute 0 -Second 0

$startRow = 37

# Row 1: Spez 1 - Rotwein - 7
$ws.Range("A$startRow").Value = $newDate
$ws.Range("B$startRow").Value = "Spez 1"
$ws.Range("C$startRow").Value = "Rotwein"
$ws.Range("D$startRow").Value = 7

# Row 2: Spez 2 - Weisswein - 7
$r = $startRow + 1
$ws.Range("A$r").Value = $newDate
$ws.Range("B$r").Value = "Spez 2"
$ws.Range("C$r").Value = "Weisswein"
$ws.Range("D$r").Value = 7

# Row 3: Spez 3 - Buch: Buchspazierer - 23.5
$r = $startRow + 2
$ws.Range("A$r").Value = $newDate
$ws.Range("B$r").Value = "Spez 3"
$ws.Range("C$r").Value = "Buch: Buchspazierer"
$ws.Range("D$r").Value = 23.5

# Resize the table / autofilter to include the new rows
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:E39"))

$ws.Range("A40").Select()
